$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 396; this shifts the existing rows 396:482
# down to 398:484 (and extends the sheet dimension to R484).
$ws.Rows.Item(396).Resize(2).Insert()

# Fill in new row 396 with the new price record
$ws.Range('A396').Value = 10
$ws.Range('B396').Value = 'Vega Modelo de Temuco'
$ws.Range('C396').Value = 'La Araucanía'
$ws.Range('D396').Value = 44889
$ws.Range('E396').Value = 9
$ws.Range('F396').Value = 100112037
$ws.Range('G396').Value = 'Cebollín'
$ws.Range('H396').Value = 'Sin especificar'
$ws.Range('I396').Value = 'Primera'
$ws.Range('J396').Value = 125
$ws.Range('K396').Value = 10000
$ws.Range('L396').Value = 11000
$ws.Range('M396').Value = 10560
$ws.Range('N396').Value = '$/docena de paquetes'
$ws.Range('O396').Value = 'Provincia de Cautín'
$ws.Range('P396').Value = 880
$ws.Range('Q396').Value = 12
$ws.Range('R396').Value = 'Hortaliza'

# Fill in new row 397 with the new price record
$ws.Range('A397').Value = 10
$ws.Range('B397').Value = 'Vega Modelo de Temuco'
$ws.Range('C397').Value = 'La Araucanía'
$ws.Range('D397').Value = 44889
$ws.Range('E397').Value = 9
$ws.Range('F397').Value = 100112037
$ws.Range('G397').Value = 'Cebollín'
$ws.Range('H397').Value = 'Sin especificar'
$ws.Range('I397').Value = 'Primera'
$ws.Range('J397').Value = 95
$ws.Range('K397').Value = 8000
$ws.Range('L397').Value = 8000
$ws.Range('M397').Value = 8000
$ws.Range('N397').Value = '$/docena de paquetes'
$ws.Range('O397').Value = "Región de O'Higgins"
$ws.Range('P397').Value = 667
$ws.Range('Q397').Value = 12
$ws.Range('R397').Value = 'Hortaliza'

# Make sure the date number format on column D carries over to the two new rows
$ws.Range('D396:D397').NumberFormat = $ws.Range('D398').NumberFormat
